$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 (shifts existing rows 80-123 down to 81-124)
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new weekly data point
$ws.Cells.Item(80, 1).Value = 7
$ws.Cells.Item(80, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(80, 3).Value = "Ñuble"
$ws.Cells.Item(80, 4).Value = 45016
$ws.Cells.Item(80, 5).Value = 16
$ws.Cells.Item(80, 6).Value = 100112037
$ws.Cells.Item(80, 7).Value = "Cebollín"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 60
$ws.Cells.Item(80, 11).Value = 7000
$ws.Cells.Item(80, 12).Value = 7000
$ws.Cells.Item(80, 13).Value = 7000
$ws.Cells.Item(80, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(80, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(80, 16).Value = 194
$ws.Cells.Item(80, 17).Value = 36
$ws.Cells.Item(80, 18).Value = "Hortaliza"
